$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns.Item(1).ColumnWidth = 14.666666666666666
$ws.Columns.Item(2).ColumnWidth = 14.666666666666666

$ws.Cells.Item(1, 1).Value = -0.12818632190988666
$ws.Cells.Item(1, 2).Value = 0.12811632393731287
$ws.Cells.Item(2, 1).Value = -0.10600997596340456
$ws.Cells.Item(2, 2).Value = 0.10580480075489884
$ws.Cells.Item(3, 1).Value = -0.073090398803492107
$ws.Cells.Item(3, 2).Value = 0.072945642980647207
$ws.Cells.Item(4, 1).Value = -0.064945642997907171
$ws.Cells.Item(4, 2).Value = 0.064408285411474253
$ws.Cells.Item(5, 1).Value = -0.061408285420397668
$ws.Cells.Item(5, 2).Value = 0.059572227854674864
$ws.Cells.Item(6, 1).Value = -0.035181288363704155
$ws.Cells.Item(6, 2).Value = 0.034748407738813825
$ws.Cells.Item(7, 1).Value = -0.024748407762665625
$ws.Cells.Item(7, 2).Value = 0.024646946890392485
$ws.Cells.Item(8, 1).Value = -0.014646946914812276
$ws.Cells.Item(8, 2).Value = 0.014477483758255971
$ws.Cells.Item(9, 1).Value = -0.012477483768626563
$ws.Cells.Item(9, 2).Value = 0.012342507374074163
$ws.Cells.Item(10, 1).Value = -0.010342507384919486
$ws.Cells.Item(10, 2).Value = 0.010333705941601323
$ws.Cells.Item(11, 1).Value = -0.013537530512276241
$ws.Cells.Item(11, 2).Value = 0.013516096404158517
$ws.Cells.Item(12, 1).Value = -0.010016096417910791
$ws.Cells.Item(12, 2).Value = 0.009863374797406177
$ws.Cells.Item(13, 1).Value = -0.0063633748116931343
$ws.Cells.Item(13, 2).Value = 0.0062952705209236015
$ws.Cells.Item(14, 1).Value = 0.0017047294564305204
$ws.Cells.Item(14, 2).Value = -0.0017269342220940231
$ws.Cells.Item(15, 1).Value = -0.0080535678444704217
$ws.Cells.Item(15, 2).Value = 0.0080347452209048598
$ws.Cells.Item(16, 1).Value = -0.0060347452327667028
$ws.Cells.Item(16, 2).Value = 0.0060037221338169466
$ws.Cells.Item(17, 1).Value = -0.0040037221459643391
$ws.Cells.Item(17, 2).Value = 0.0039999999841100475
$ws.Cells.Item(18, 1).Value = -0.016106347985207492
$ws.Cells.Item(18, 2).Value = 0.016091868351406191
$ws.Cells.Item(19, 1).Value = -0.012091868359050739
$ws.Cells.Item(19, 2).Value = 0.012017182703600415
$ws.Cells.Item(20, 1).Value = -0.0080171827118373784
$ws.Cells.Item(20, 2).Value = 0.0080057087990166309
$ws.Cells.Item(21, 1).Value = -0.0040057088073286451
$ws.Cells.Item(21, 2).Value = 0.003999999991604497
$ws.Cells.Item(22, 1).Value = -0.028714401960609592
$ws.Cells.Item(22, 2).Value = 0.028431125512870636
$ws.Cells.Item(23, 1).Value = -0.040493416211369926
$ws.Cells.Item(23, 2).Value = 0.040098003022554174
$ws.Cells.Item(24, 1).Value = -0.020098003062550163
$ws.Cells.Item(24, 2).Value = 0.019999999959466663
$ws.Cells.Item(25, 1).Value = -0.022390939499958407
$ws.Cells.Item(25, 2).Value = 0.022337136006626679
$ws.Cells.Item(26, 1).Value = -0.019837136016624513
$ws.Cells.Item(26, 2).Value = 0.019771881668303237
$ws.Cells.Item(27, 1).Value = -0.01727188167849647
$ws.Cells.Item(27, 2).Value = 0.016910737195145131
$ws.Cells.Item(28, 1).Value = -0.014910737205113378
$ws.Cells.Item(28, 2).Value = 0.014682629219992194
$ws.Cells.Item(29, 1).Value = -0.0076826292395315576
$ws.Cells.Item(29, 2).Value = 0.0076271685843254389
$ws.Cells.Item(30, 1).Value = 0.052372831299137612
$ws.Cells.Item(30, 2).Value = -0.052634545024604407
$ws.Cells.Item(31, 1).Value = 0.059634545006710837
$ws.Cells.Item(31, 2).Value = -0.059730942542907783
$ws.Cells.Item(32, 1).Value = -0.0040015814419742668
$ws.Cells.Item(32, 2).Value = 0.0039999999880784287
